$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'36.517.91"
$ws.Range("E2").Value = "'  -1.22%  "
$ws.Range("D3").Value = "'2.056.05"
$ws.Range("E4").Value = "'  +0.07%  "
$ws.Range("D5").Value = "'242.63"
$ws.Range("E5").Value = "'  -1.13%  "
$ws.Range("E6").Value = "'  +0.92%  "
$ws.Range("E7").Value = "'  +0.03%  "
$ws.Range("D8").Value = "'54.69"
$ws.Range("E8").Value = "'  -6.09%  "
$ws.Range("D9").Value = "'58.44"
$ws.Range("E9").Value = "'  -2.29%  "
$ws.Range("E10").Value = "'  -4.30%  "
$ws.Range("E11").Value = "'  -2.40%  "
$ws.Range("E12").Value = "'  -2.92%  "
$ws.Range("D13").Value = "'0.912"
$ws.Range("E13").Value = "'  +3.43%  "
$ws.Range("D14").Value = "'14.75"
$ws.Range("E14").Value = "'  -4.79%  "
$ws.Range("D15").Value = "'2.356.38"
$ws.Range("E15").Value = "'  +0.65%  "
$ws.Range("D16").Value = "'5.41"
$ws.Range("E16").Value = "'  -4.43%  "
$ws.Range("D17").Value = "'2.027.54"
$ws.Range("E17").Value = "'  -0.67%  "
$ws.Range("D18").Value = "'36.468.61"
$ws.Range("E18").Value = "'  -1.30%  "
$ws.Range("D19").Value = "'16.77"
$ws.Range("E19").Value = "'  -7.83%  "
$ws.Range("D20").Value = "'71.96"
$ws.Range("E20").Value = "'  -2.47%  "
$ws.Range("E21").Value = "'  -3.53%  "
$ws.Range("D22").Value = "'238.52"
$ws.Range("E22").Value = "'  +1.10%  "
$ws.Range("D23").Value = "'5.25"
$ws.Range("E23").Value = "'  -2.66%  "
$ws.Range("E24").Value = "'  +0.16%  "
$ws.Range("D25").Value = "'2.36"
$ws.Range("E25").Value = "'  -3.99%  "
$ws.Range("E26").Value = "'  -2.79%  "
$ws.Range("E27").Value = "'  +0.15%  "
$ws.Range("D28").Value = "'164.81"
$ws.Range("E28").Value = "'  -2.83%  "
$ws.Range("D29").Value = "'20.09"
$ws.Range("E29").Value = "'  +0.77%  "
$ws.Range("B30").Value = "'ImmutableX"
$ws.Range("C30").Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D30").Value = "'1.23"
$ws.Range("E30").Value = "'  +11.29%  "
$ws.Range("B31").Value = "'Stellar"
$ws.Range("C31").Value = "'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D31").Value = "'0.122"
$ws.Range("E31").Value = "'  -1.09%  "
$ws.Range("D32").Value = "'5.10"
$ws.Range("E32").Value = "'  -4.99%  "
$ws.Range("D33").Value = "'4.46"
$ws.Range("E33").Value = "'  -5.22%  "
$ws.Range("D34").Value = "'0.0595"
$ws.Range("E34").Value = "'  -2.81%  "
$ws.Range("E35").Value = "'  +0.14%  "
$ws.Range("E36").Value = "'  -0.68%  "
$ws.Range("E37").Value = "'  -1.28%  "
$ws.Range("D38").Value = "'0.0820"
$ws.Range("E38").Value = "'  -5.96%  "
$ws.Range("D39").Value = "'1.25"
$ws.Range("E39").Value = "'  -5.11%  "
$ws.Range("D40").Value = "'4.84"
$ws.Range("E40").Value = "'  -4.89%  "
$ws.Range("E41").Value = "'  -2.82%  "
$ws.Range("B42").Value = "'HuobiToken"
$ws.Range("C42").Value = "'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D42").Value = "'2.84"
$ws.Range("E42").Value = "'  -8.73%  "
$ws.Range("B43").Value = "'Cronos"
$ws.Range("C43").Value = "'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D43").Value = "'0.0929"
$ws.Range("E43").Value = "'  -5.45%  "
$ws.Range("D44").Value = "'1.11"
$ws.Range("E44").Value = "'  -2.77%  "
$ws.Range("D45").Value = "'93.88"
$ws.Range("E45").Value = "'  -3.18%  "
$ws.Range("D46").Value = "'1.414.67"
$ws.Range("E46").Value = "'  +9.29%  "
$ws.Range("D47").Value = "'7.62"
$ws.Range("E47").Value = "'  +13.24%  "
$ws.Range("D48").Value = "'15.95"
$ws.Range("E48").Value = "'  -5.41%  "
$ws.Range("E49").Value = "'  +0.06%  "
$ws.Range("D50").Value = "'2.27"
$ws.Range("E50").Value = "'  -2.57%  "
$ws.Range("D51").Value = "'2.241.26"
$ws.Range("E51").Value = "'  +0.71%  "
